$d = $word.ActiveDocument

# 1. Table width: auto -> 100% (pct, w=5000)
#    PreferredWidth (for percent type) maps to OOXML w:w = PreferredWidth * 20,
#    so 250 -> w:w="5000" (5000 fiftieths-of-a-percent == 100%).
$t = $d.Tables.Item(1)
$t.PreferredWidthType = 2   # wdPreferredWidthPercent
$t.PreferredWidth = 250

# Namespace declaration needed for InsertXML fragments
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# 2. Fill in the three rubric cells on the second (body) row, left-justifying
#    the paragraph and inserting the rubric description text.
$row = 2

$cell1 = $t.Cell($row, 1)
[void]$cell1.Range.InsertXML("<w:p $wNs><w:pPr><w:pStyle w:val='Compact'/><w:jc w:val='left'/></w:pPr><w:r><w:t xml:space='preserve'>Response directly addresses the prompt with specific details from the readings and NotebookLM. Includes concrete examples from field experience or teaching practice.</w:t></w:r></w:p>")

$cell2 = $t.Cell($row, 2)
[void]$cell2.Range.InsertXML("<w:p $wNs><w:pPr><w:pStyle w:val='Compact'/><w:jc w:val='left'/></w:pPr><w:r><w:t xml:space='preserve'>Response addresses the prompt but lacks specific details or examples. May be vague or general.</w:t></w:r></w:p>")

$cell3 = $t.Cell($row, 3)
[void]$cell3.Range.InsertXML("<w:p $wNs><w:pPr><w:pStyle w:val='Compact'/><w:jc w:val='left'/></w:pPr><w:r><w:t xml:space='preserve'>No response or response does not address the prompt.</w:t></w:r></w:p>")
